$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as text in the source data (e.g. "65.830.56", "1.00"),
# so a leading apostrophe forces Excel to keep them as text instead of parsing as numbers.

$ws.Range("D2").Value = "'65.830.56"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "'3.306.34"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'185.73"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'554.56"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'3.303.05"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("E9").Value = "  -2.75%  "
$ws.Range("E10").Value = "  -5.56%  "
$ws.Range("D11").Value = "'0.577"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").Value = "'45.69"
$ws.Range("E12").Value = "  -3.23%  "
$ws.Range("D13").Value = "'0.0000263"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "'3.830.72"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").Value = "'8.44"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "'574.00"
$ws.Range("E16").Value = "  -8.18%  "
$ws.Range("D17").Value = "'65.853.80"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "'3.297.87"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "'17.66"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  -4.49%  "
$ws.Range("D22").Value = "'0.889"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").Value = "'18.02"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").Value = "'97.67"
$ws.Range("E25").Value = "  -7.85%  "
$ws.Range("D26").Value = "'3.93"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").Value = "'9.35"
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'8.39"
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'30.43"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("E31").Value = "  +5.49%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").Value = "'3.68"
$ws.Range("E32").Value = "  -8.54%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'563.14"
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("D34").Value = "'10.80"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("D36").Value = "'3.719.58"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'55.44"
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").Value = "'33.76"
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("D41").Value = "'0.0₃0684"
$ws.Range("E41").Value = "  -5.83%  "
$ws.Range("E42").Value = "  -7.85%  "
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("E44").Value = "  -4.66%  "
$ws.Range("D45").Value = "'0.333"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("B49").Value = "CoreDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D49").Value = "'2.91"
$ws.Range("E49").Value = "  -13.25%  "
$ws.Range("D50").Value = "'2.51"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").Value = "'127.07"
$ws.Range("E51").Value = "  +4.43%  "
